# Upload stage 3 evidence for B1-B2
#
# Narrative: the user was on the "Info" sheet (cell B14 selected), clicked
# into Info!B2, then went to sheet "B1" and filled in the two evidence
# TxHash cells (A2/A3), leaving the selection on A4 just past the data.
# Finally they moved to sheet "B2", filled in its A2/A3 TxHash cells, and
# left the selection on J12 - this is the sheet/cell that was active when
# the workbook was saved.

$wb = $excel.ActiveWorkbook

# Start: touch the Info sheet selection (B14 -> B2) before navigating away.
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Activate()
$wsInfo.Range("B2").Select()

# Fill in stage-3 evidence on sheet "B1".
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Activate()
$wsB1.Range("A2").Value = "E5AD0CE00693DB7C9226CB84BC5E3CC5755CE283C76712AD412EA4F5C91542E8"
$wsB1.Range("A3").Value = "03F24C1272F863EA37124400DE343E3C64263B822A38ED597443C7FAB86BA048"
$wsB1.Range("A4").Select()

# Fill in stage-3 evidence on sheet "B2" (the sheet left active on save).
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Activate()
$wsB2.Range("A2").Value = "400486D619EE56A7068F348F8EF0BF80EEA9C665BBBAFF2BDBB68217E71EFE51"
$wsB2.Range("A3").Value = "486586F1DEB8FE28694CB4FE0AEF4A80CA17F4BE537359E82440EFB277B6E741"
$wsB2.Range("J12").Select()
